$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.839.76'
$ws.Range("E2").Value = '  +0.41%  '
$ws.Range("D3").Value = '2.437.91'
$ws.Range("E3").Value = '  -0.90%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '559.56'
$ws.Range("E5").Value = '  +0.20%  '
$ws.Range("D6").Value = '162.08'
$ws.Range("E6").Value = '  +0.80%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = '0.513'
$ws.Range("E8").Value = '  +1.37%  '
$ws.Range("E9").Value = '  +11.87%  '
$ws.Range("E10").Value = '  -1.57%  '
$ws.Range("E11").Value = '  +0.36%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.60'
$ws.Range("E12").Value = '  -5.08%  '
$ws.Range("E13").Value = '  +4.97%  '
$ws.Range("D14").Value = '68.723.07'
$ws.Range("E14").Value = '  +0.34%  '
$ws.Range("D15").Value = '2.886.09'
$ws.Range("E15").Value = '  -0.68%  '
$ws.Range("D16").Value = '23.17'
$ws.Range("E16").Value = '  -0.59%  '
$ws.Range("D17").Value = '2.437.65'
$ws.Range("E17").Value = '  -1.64%  '
$ws.Range("D18").Value = '10.53'
$ws.Range("E18").Value = '  -0.75%  '
$ws.Range("D19").Value = '338.56'
$ws.Range("E19").Value = '  +1.50%  '
$ws.Range("D20").Value = '6.93'
$ws.Range("E20").Value = '  +0.55%  '
$ws.Range("E21").Value = '  +1.55%  '
$ws.Range("D22").Value = '1.92'
$ws.Range("E22").Value = '  +2.71%  '
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("E24").Value = '  +1.00%  '
$ws.Range("D25").Value = '3.71'
$ws.Range("E25").Value = '  +2.33%  '
$ws.Range("D26").Value = '2.564.74'
$ws.Range("D27").Value = '1.01'
$ws.Range("E27").Value = '  +0.61%  '
$ws.Range("D28").Value = '8.17'
$ws.Range("E28").Value = '  +0.51%  '
$ws.Range("D29").Value = '0.0₃0818'
$ws.Range("E29").Value = '  +0.58%  '
$ws.Range("D30").Value = '7.13'
$ws.Range("E30").Value = '  -0.63%  '
$ws.Range("E31").Value = '  +0.04%  '
$ws.Range("D32").Value = '427.79'
$ws.Range("E32").Value = '  +0.21%  '
$ws.Range("E33").Value = '  +2.12%  '
$ws.Range("E34").Value = '  -0.12%  '
$ws.Range("D35").Value = '160.29'
$ws.Range("E35").Value = '  +1.70%  '
$ws.Range("D36").Value = '19.01'
$ws.Range("E36").Value = '  +0.11%  '
$ws.Range("E37").Value = '  +0.01%  '
$ws.Range("D38").Value = '17.96'
$ws.Range("E38").Value = '  +1.46%  '
$ws.Range("E39").Value = '  -2.33%  '
$ws.Range("E40").Value = '  -0.60%  '
$ws.Range("E41").Value = '  +3.35%  '
$ws.Range("E42").Value = '  -1.13%  '
$ws.Range("D43").Value = '1.07'
$ws.Range("E43").Value = '  +0.59%  '
$ws.Range("E44").Value = '  -0.36%  '
$ws.Range("E45").Value = '  +0.23%  '
$ws.Range("D46").Value = '129.94'
$ws.Range("E46").Value = '  +0.78%  '
$ws.Range("E47").Value = '  +0.58%  '
$ws.Range("E48").Value = '  +0.15%  '
$ws.Range("E49").Value = '  -0.45%  '
$ws.Range("D50").Value = '0.0924'
$ws.Range("E50").Value = '  +1.76%  '
$ws.Range("E51").Value = '  +1.07%  '
